# Parameter-editing script:
#   - Row 15 ("qViewTrialByTrial") is repurposed to hold the new
#     "max_onset_allowed (N_sample)" group parameter (default 1, range 0..10000).
#   - Three new rows are appended for the remaining new group parameters,
#     with "qViewTrialByTrial" moved down to become the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 15: was qViewTrialByTrial, now max_onset_allowed ---
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Name"
$ws.Cells.Item(15, 3).Value = "Group data proc: `nmax_onset_allowed (N_sample)"
$ws.Cells.Item(15, 4).Value = "RangeInputBox"
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 10000

# --- New row 16: min_diff_allowed (second) ---
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Name"
$ws.Cells.Item(16, 3).Value = "Group data proc: `nmin_diff_allowed (second)"
$ws.Cells.Item(16, 4).Value = "RangeInputBox"
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 3

# --- New row 17: max_diff_allowed (second) ---
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "Name"
$ws.Cells.Item(17, 3).Value = "Group data proc: `nmax_diff_allowed (second)"
$ws.Cells.Item(17, 4).Value = "RangeInputBox"
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 3

# --- New row 18: qViewTrialByTrial (moved down from old row 15) ---
$ws.Cells.Item(18, 1).Value = 0
$ws.Cells.Item(18, 2).Value = "Name"
$ws.Cells.Item(18, 3).Value = "Group data proc: `nqViewTrialByTrial"
$ws.Cells.Item(18, 4).Value = "RangeInputBox"
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 1

# Keep row heights on their default (un-customized) auto height, same as
# every other row in the sheet, rather than leaving an explicit wrapped
# row height behind after writing the multi-line labels above.
$ws.Rows("15:18").AutoFit() | Out-Null
